$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D column (Price) to Text format first so numeric-looking
# strings (e.g. "225.90") are not auto-converted to floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "34.148.82"
$ws.Range("E2").Value = "  +0.87%  "

# Row 3
$ws.Range("D3").Value = "1.777.70"
$ws.Range("E3").Value = "  -0.10%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").Value = "225.90"
$ws.Range("E5").Value = "  +0.73%  "

# Row 6
$ws.Range("D6").Value = "0.546"
$ws.Range("E6").Value = "  +0.11%  "

# Row 7
$ws.Range("E7").Value = "  +0.18%  "

# Row 8
$ws.Range("D8").Value = "31.57"
$ws.Range("E8").Value = "  -0.13%  "

# Row 9
$ws.Range("D9").Value = "0.291"
$ws.Range("E9").Value = "  +0.87%  "

# Row 10
$ws.Range("D10").Value = "0.0691"
$ws.Range("E10").Value = "  +2.05%  "

# Row 11
$ws.Range("E11").Value = "  +1.09%  "

# Row 12
$ws.Range("D12").Value = "2.036.13"
$ws.Range("E12").Value = "  +0.25%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.784.08"
$ws.Range("E13").Value = "  +0.10%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "10.91"
$ws.Range("E14").Value = "  -1.93%  "

# Row 15
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "34.134.79"
$ws.Range("E15").Value = "  +0.82%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "0.621"
$ws.Range("E16").Value = "  +2.09%  "

# Row 17
$ws.Range("D17").Value = "4.18"
$ws.Range("E17").Value = "  +1.13%  "

# Row 18
$ws.Range("D18").Value = "67.80"
$ws.Range("E18").Value = "  +1.56%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0799"
$ws.Range("E19").Value = "  +3.82%  "

# Row 20
$ws.Range("D20").Value = "245.35"
$ws.Range("E20").Value = "  +2.73%  "

# Row 21
$ws.Range("D21").Value = "10.97"
$ws.Range("E21").Value = "  +4.12%  "

# Row 22
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("E23").Value = "  +1.92%  "

# Row 24
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  -1.35%  "

# Row 25
$ws.Range("D25").Value = "162.53"
$ws.Range("E25").Value = "  +0.94%  "

# Row 26
$ws.Range("E26").Value = "  +2.37%  "

# Row 27
$ws.Range("D27").Value = "16.28"
$ws.Range("E27").Value = "  +1.25%  "

# Row 28
$ws.Range("E28").Value = "  +1.87%  "

# Row 29
$ws.Range("E29").Value = "  +0.18%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.23"
$ws.Range("E30").Value = "  +0.85%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.0520"
$ws.Range("E31").Value = "  +2.19%  "

# Row 32
$ws.Range("D32").Value = "3.73"
$ws.Range("E32").Value = "  +4.15%  "

# Row 33
$ws.Range("D33").Value = "3.70"
$ws.Range("E33").Value = "  +5.04%  "

# Row 34
$ws.Range("D34").Value = "1.79"
$ws.Range("E34").Value = "  -1.14%  "

# Row 35
$ws.Range("D35").Value = "1.437.88"
$ws.Range("E35").Value = "  +3.56%  "

# Row 36
$ws.Range("D36").Value = "0.657"
$ws.Range("E36").Value = "  +3.32%  "

# Row 37
$ws.Range("E37").Value = "  +6.53%  "

# Row 38
$ws.Range("D38").Value = "0.0190"
$ws.Range("E38").Value = "  +2.81%  "

# Row 39
$ws.Range("D39").Value = "1.05"
$ws.Range("E39").Value = "  +0.44%  "

# Row 40
$ws.Range("D40").Value = "80.08"
$ws.Range("E40").Value = "  +2.32%  "

# Row 41
$ws.Range("E41").Value = "  -0.25%  "

# Row 42
$ws.Range("D42").Value = "0.920"
$ws.Range("E42").Value = "  +0.84%  "

# Row 43
$ws.Range("D43").Value = "2.66"
$ws.Range("E43").Value = "  +0.69%  "

# Row 44
$ws.Range("D44").Value = "13.47"
$ws.Range("E44").Value = "  -0.05%  "

# Row 45
$ws.Range("D45").Value = "0.0510"
$ws.Range("E45").Value = "  -0.24%  "

# Row 46
$ws.Range("E46").Value = "  +3.83%  "

# Row 47
$ws.Range("E47").Value = "  +0.11%  "

# Row 48
$ws.Range("D48").Value = "0.0₆0135"
$ws.Range("E48").Value = "  -0.09%  "

# Row 49
$ws.Range("D49").Value = "1.937.80"
$ws.Range("E49").Value = "  +0.08%  "

# Row 50
$ws.Range("D50").Value = "104.19"
$ws.Range("E50").Value = "  -0.94%  "

# Row 51
$ws.Range("E51").Value = "  +0.18%  "

# Restore normal cell style on the Price column (keeps values as text,
# drops the temporary Text number-format style so cell styling matches original).
$ws.Range("D2:D51").Style = "Normal"

Write-Output "All updates applied."
